$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the nowcast values for the existing rows (2-11, columns B-K)
# with the latest model run. Row 1 headers and the column-A report dates
# for these rows are unchanged.
$ws.Range("B2").Value = 0.30772807106524996
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("B3").Value = 0.3184713988597034
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.009726066265945011
$ws.Range("E3").Value = 0.00011512948610463593
$ws.Range("F3").Value = -0.00019105368180970554
$ws.Range("G3").Value = 0.0012638823234799462
$ws.Range("H3").Value = -0.00010930112311942647
$ws.Range("I3").Value = -0.0001645291645544686
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.00010313368840741965
$ws.Range("B4").Value = 0.2949145877630909
$ws.Range("C4").Value = -0.0073884054735758355
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.00008506133831743608
$ws.Range("F4").Value = -0.0003765189589069301
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = -0.0002907573987097591
$ws.Range("I4").Value = -0.014164533544282477
$ws.Range("J4").Value = -0.000977968703648797
$ws.Range("K4").Value = -0.0002735656791712504
$ws.Range("B5").Value = 0.3332019896340785
$ws.Range("C5").Value = 0.04337482203201637
$ws.Range("D5").Value = 0.0012266331878305052
$ws.Range("E5").Value = 0.00005465937491974757
$ws.Range("F5").Value = 0.00034998010258263716
$ws.Range("G5").Value = -0.0015671520173489756
$ws.Range("H5").Value = -0.0003948354263689205
$ws.Range("I5").Value = -0.004916961536850888
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.0001602561542071168
$ws.Range("B6").Value = 0.6173441306480913
$ws.Range("C6").Value = 0.3075698053168663
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.00010167181955926962
$ws.Range("F6").Value = 0.0007718350562329246
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.00012767888887290998
$ws.Range("I6").Value = -0.027089523327448052
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.0029160310376753484
$ws.Range("B7").Value = 0.555468924472122
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.06948820452310353
$ws.Range("E7").Value = -0.000020767087092853857
$ws.Range("F7").Value = -0.0033255933602265434
$ws.Range("G7").Value = 0.008546282272349678
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.0030916082536977365
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.0006785317315938233
$ws.Range("B8").Value = 0.5666348005563879
$ws.Range("C8").Value = 0.05082618773676498
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.00019638827011087863
$ws.Range("F8").Value = -0.012385753190251597
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = -0.0007726167823217891
$ws.Range("I8").Value = -0.034340517716496655
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.007642187766460085
$ws.Range("B9").Value = 0.11597676924464881
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = -0.30994738914408965
$ws.Range("E9").Value = -0.005438166767006645
$ws.Range("F9").Value = -0.13491102601873015
$ws.Range("G9").Value = -0.011553636889168752
$ws.Range("H9").Value = -0.0021166867756976706
$ws.Range("I9").Value = 0.012722945646474356
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0.000585928636479438
$ws.Range("B10").Value = 0.13361113710843997
$ws.Range("C10").Value = 0.01992998162094068
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -0.0004532487040947844
$ws.Range("F10").Value = 0.0015372809472543793
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.0003305213294890289
$ws.Range("I10").Value = 0.02108615174898567
$ws.Range("J10").Value = -0.02314626654877869
$ws.Range("K10").Value = -0.0016500525300051222
$ws.Range("B11").Value = 0.44477962946062544
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.2761583840906835
$ws.Range("E11").Value = -0.006043481514773706
$ws.Range("F11").Value = 0.014954540765309427
$ws.Range("G11").Value = -0.006731961932654688
$ws.Range("H11").Value = -0.006189389416026967
$ws.Range("I11").Value = 0.09066576323380128
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = -0.05164536287415339

# Append the new row for the latest report date (2025-08-30).
# Pre-format column A as text so the ISO date string is not
# auto-converted to a date serial, then drop back to the default
# (General) style so the new row matches the rest of the table.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-08-30"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = 0.21532128703895853
$ws.Range("C12").Value = -0.16089031005936275
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.0008385778461936903
$ws.Range("F12").Value = -0.0019084763498543732
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = -0.00198811206248848
$ws.Range("I12").Value = -0.046208053274959734
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -0.019301968521195256

# Column widths widened slightly to fit the refreshed values
$ws.Columns.Item(3).ColumnWidth = 15.24609375
$ws.Columns.Item(5).ColumnWidth = 16.24609375
